# The commit swaps the contents of ppt/theme/theme1.xml (the slide
# master's theme, originally the "Integral" color theme) and
# ppt/theme/theme2.xml (the notes master's theme, originally the
# default "Office Theme" color theme) - i.e. the slide deck's design
# ends up using the stock "Office Theme" palette instead of "Integral".
#
# Drive this the same way a user would from the Design tab: recolour
# the presentation's theme (ThemeColorScheme) to the 12 standard
# "Office" theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink). PowerPoint's RGB() packs colours as R + G*256 + B*65536,
# i.e. the same value you'd get from RGB(r,g,b) for each srgbClr below:
#   1  dk1       000000 -> 0
#   2  lt1       FFFFFF -> 16777215
#   3  dk2       44546A -> 6968388
#   4  lt2       E7E6E6 -> 15132391
#   5  accent1   5B9BD5 -> 13998939
#   6  accent2   ED7D31 -> 3243501
#   7  accent3   A5A5A5 -> 10855845
#   8  accent4   FFC000 -> 49407
#   9  accent5   4472C4 -> 12874308
#   10 accent6   70AD47 -> 4697456
#   11 hlink     0563C1 -> 12673797
#   12 folHlink  954F72 -> 7491477

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB = 0
$tcs.Item(2).RGB = 16777215
$tcs.Item(3).RGB = 6968388
$tcs.Item(4).RGB = 15132391
$tcs.Item(5).RGB = 13998939
$tcs.Item(6).RGB = 3243501
$tcs.Item(7).RGB = 10855845
$tcs.Item(8).RGB = 49407
$tcs.Item(9).RGB = 12874308
$tcs.Item(10).RGB = 4697456
$tcs.Item(11).RGB = 12673797
$tcs.Item(12).RGB = 7491477
